$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.631.89"
$ws.Range("E2").Value = "  +3.23%  "
$ws.Range("D3").Value = "2.441.19"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.08%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "2.439.82"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.91%  "
$ws.Range("E15").Value = "  +5.15%  "
$ws.Range("D16").Value = "2.883.27"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("D17").Value = "62.469.32"
$ws.Range("E17").Value = "  +3.29%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.441.28"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("B19").Value = "BabyDogeCoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D19").Value = "0.0₆0950"
$ws.Range("E19").Value = "  +231.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "325.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("E24").Value = "  +9.98%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "648.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.99%  "
$ws.Range("E28").Value = "  +13.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.84%  "
$ws.Range("D30").Value = "0.0₃0979"
$ws.Range("E30").Value = "  +4.37%  "
$ws.Range("D31").Value = "2.558.48"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.76%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.139"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.43%  "
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.18%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "152.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.372"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("E43").Value = "  +8.42%  "
$ws.Range("E44").Value = "  +4.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +28.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "144.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.602"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.93%  "
